$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values in row 2
$ws.Range("D2").Value = 5
$ws.Range("F2").Value = 3
$ws.Range("H2").Value = 46

# Update the active selection to D2
$ws.Range("D2").Select() | Out-Null
